# Weekly update: a new price observation was recorded for the week of
# 2021-10-12 ("Región de O'Higgins"). It is inserted as the new row 11
# (right after the fixed preamble rows 1-10), pushing every existing
# data row (old 11..37) down by one (new 12..38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11:37 down to 12:38, leaving a blank row 11 to fill in.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "Macroferia Regional de Talca"
$ws.Range("C11").Value = "Maule"
$ws.Range("D11").Value = (Get-Date -Year 2021 -Month 10 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 100112022
$ws.Range("G11").Value = "Arveja Verde"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 23000
$ws.Range("L11").Value = 23000
$ws.Range("M11").Value = 23000
$ws.Range("N11").Value = "`$/saco 25 kilos"
$ws.Range("O11").Value = "Región de O'Higgins"
$ws.Range("P11").Value = 920
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
